$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (shifts existing rows 8..108 down to 9..109)
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new record
$ws.Range("A8").Value = 1
$ws.Range("B8").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C8").Value = "Arica y Parinacota"
$ws.Range("D8").Value = 44901
$ws.Range("E8").Value = 15
$ws.Range("F8").Value = 100114001
$ws.Range("G8").Value = "Papa"
$ws.Range("H8").Value = "Red Lady"
$ws.Range("I8").Value = "1a nueva(o)"
$ws.Range("J8").Value = 1000
$ws.Range("K8").Value = 18000
$ws.Range("L8").Value = 20000
$ws.Range("M8").Value = 19000
$ws.Range("N8").Value = "`$/saco 25 kilos"
$ws.Range("O8").Value = "Provincia de Melipilla"
$ws.Range("P8").Value = 760
$ws.Range("Q8").Value = 25
$ws.Range("R8").Value = "Hortaliza"
